$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the formatting
# (style) already used by the existing header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I and J, rows 2-30
$data = @{
    2  = @(6, 8)
    3  = @(4, 6)
    4  = @(6, 7)
    5  = @(2, 3)
    6  = @(7, 7)
    7  = @(5, 7)
    8  = @(11, 14)
    9  = @(7, 9)
    10 = @(9, 9)
    11 = @(4, 7)
    12 = @(10, 10)
    13 = @(6, 7)
    14 = @(6, 7)
    15 = @(7, 9)
    16 = @(6, 8)
    17 = @(5, 8)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 5)
    21 = @(1, 6)
    22 = @(1, 5)
    23 = @(1, 5)
    24 = @(1, 4)
    25 = @(1, 6)
    26 = @(1, 5)
    27 = @(1, 4)
    28 = @(1, 4)
    29 = @(1, 3)
    30 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
